$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AutoIncrement")

# Bump the running sequence number that drives every "AB-03" / "ABs-03" /
# "SGTTAP-VNTTVN-AB-03-003" style generated id across the workbook.
$ws.Range("A2").Value = "04"

# Leave the sheet's selection where the author left it after making the edit.
$ws.Range("A3").Select()
